$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 964, shifting existing rows 964:1071 down to 965:1072.
$ws.Rows("964:964").Insert()

# Populate the newly inserted row 964 with the new data point.
$ws.Range("A964").Value = 7
$ws.Range("B964").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C964").Value = "Ñuble"
$ws.Range("D964").Value = 45194
$ws.Range("E964").Value = 16
$ws.Range("F964").Value = "Fruta"
$ws.Range("G964").Value = 100108
$ws.Range("H964").Value = "Tropicales y subtropicales"
$ws.Range("I964").Value = 100108006
$ws.Range("J964").Value = "Plátano"
$ws.Range("K964").Value = "Sin especificar"
$ws.Range("L964").Value = "Pintón"
$ws.Range("M964").Value = 200
$ws.Range("N964").Value = 24000
$ws.Range("O964").Value = 24000
$ws.Range("P964").Value = 24000
$ws.Range("Q964").Value = "$/caja 20 kilos"
$ws.Range("R964").Value = "Ecuador"
$ws.Range("S964").Value = 1200
$ws.Range("T964").Value = 20
